$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (these are display strings, not numbers).
$textCells = @("D4","D5","D7","D8","D9","D10","D11","D12","D14","D16","D17","D19","D22","D25","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D41","D42","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value2 = "30.536.75"
$ws.Range("E2").Value2 = "  +0.28%  "
$ws.Range("D3").Value2 = "2.112.12"
$ws.Range("E3").Value2 = "  +0.96%  "
$ws.Range("D4").Value2 = "1.011"
$ws.Range("E4").Value2 = "  +0.75%  "
$ws.Range("D5").Value2 = "335.88"
$ws.Range("E5").Value2 = "  +1.82%  "
$ws.Range("E6").Value2 = "  +0.62%  "
$ws.Range("D7").Value2 = "0.5241"
$ws.Range("E7").Value2 = "  +0.67%  "
$ws.Range("D8").Value2 = "0.4550"
$ws.Range("E8").Value2 = "  +4.11%  "
$ws.Range("D9").Value2 = "54.97"
$ws.Range("E9").Value2 = "  +1.58%  "
$ws.Range("D10").Value2 = "0.09109"
$ws.Range("E10").Value2 = "  +3.02%  "
$ws.Range("D11").Value2 = "1.170"
$ws.Range("E11").Value2 = "  +1.57%  "
$ws.Range("D12").Value2 = "24.54"
$ws.Range("E12").Value2 = "  +1.24%  "
$ws.Range("D13").Value2 = "2.118.59"
$ws.Range("E13").Value2 = "  +2.23%  "
$ws.Range("D14").Value2 = "6.842"
$ws.Range("E14").Value2 = "  +2.26%  "
$ws.Range("E15").Value2 = "  +5.72%  "
$ws.Range("D16").Value2 = "0.00001176"
$ws.Range("E16").Value2 = "  +5.03%  "
$ws.Range("D17").Value2 = "97.02"
$ws.Range("E17").Value2 = "  +1.36%  "
$ws.Range("E18").Value2 = "  +0.61%  "
$ws.Range("D19").Value2 = "0.06682"
$ws.Range("E19").Value2 = "  +1.30%  "
$ws.Range("E20").Value2 = "  +0.85%  "
$ws.Range("D22").Value2 = "6.270"
$ws.Range("E22").Value2 = "  +0.00%  "
$ws.Range("D23").Value2 = "30.603.00"
$ws.Range("E23").Value2 = "  +0.40%  "
$ws.Range("E24").Value2 = "  +4.28%  "
$ws.Range("D25").Value2 = "2.354"
$ws.Range("E25").Value2 = "  +0.62%  "
$ws.Range("D26").Value2 = "2.365.32"
$ws.Range("E26").Value2 = "  +2.00%  "
$ws.Range("D27").Value2 = "22.30"
$ws.Range("E27").Value2 = "  +0.36%  "
$ws.Range("D28").Value2 = "163.64"
$ws.Range("E28").Value2 = "  +0.47%  "
$ws.Range("E29").Value2 = "  -1.38%  "
$ws.Range("D30").Value2 = "133.86"
$ws.Range("E30").Value2 = "  +1.72%  "
$ws.Range("D31").Value2 = "1.209"
$ws.Range("E31").Value2 = "  +2.20%  "
$ws.Range("D32").Value2 = "0.1069"
$ws.Range("E32").Value2 = "  +0.30%  "
$ws.Range("D33").Value2 = "1.641"
$ws.Range("E33").Value2 = "  -0.01%  "
$ws.Range("D34").Value2 = "6.362"
$ws.Range("E34").Value2 = "  +3.24%  "
$ws.Range("D35").Value2 = "3.948"
$ws.Range("E35").Value2 = "  +1.11%  "
$ws.Range("D36").Value2 = "10.54"
$ws.Range("E36").Value2 = "  +4.69%  "
$ws.Range("D37").Value2 = "5.902"
$ws.Range("E37").Value2 = "  +8.28%  "
$ws.Range("D38").Value2 = "0.02616"
$ws.Range("E38").Value2 = "  +1.43%  "
$ws.Range("D39").Value2 = "0.06810"
$ws.Range("E39").Value2 = "  +0.25%  "
$ws.Range("E40").Value2 = "  +3.22%  "
$ws.Range("D41").Value2 = "12.56"
$ws.Range("E41").Value2 = "  -0.58%  "
$ws.Range("D42").Value2 = "0.6862"
$ws.Range("E42").Value2 = "  -0.14%  "
$ws.Range("E43").Value2 = "  +0.07%  "
$ws.Range("D44").Value2 = "15.00"
$ws.Range("E44").Value2 = "  +7.73%  "
$ws.Range("D45").Value2 = "0.6438"
$ws.Range("E45").Value2 = "  +1.55%  "
$ws.Range("D46").Value2 = "2.304"
$ws.Range("E46").Value2 = "  +5.03%  "
$ws.Range("D47").Value2 = "3.683"
$ws.Range("E47").Value2 = "  +1.58%  "
$ws.Range("D48").Value2 = "0.00000000359"
$ws.Range("E48").Value2 = "  +20.07%  "
$ws.Range("E49").Value2 = "  +0.69%  "
$ws.Range("D50").Value2 = "83.19"
$ws.Range("E50").Value2 = "  +1.87%  "
$ws.Range("D51").Value2 = "0.3342"
$ws.Range("E51").Value2 = "  +12.23%  "
